# Apply "scale question, showing filled tests" update to WorkReport.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data (Task description, Hours, Date)
$newRows = @(
    @{ Row = 45; Text = "Setting up database project, writing what to do, overview"; Hours = 3;   Date = (Get-Date -Year 2012 -Month 2 -Day 16 -Hour 0 -Minute 0 -Second 0) },
    @{ Row = 46; Text = "Scale question gui and controls";                           Hours = 4;   Date = (Get-Date -Year 2012 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0) },
    @{ Row = 47; Text = "Scale question, showing filled tests list and window";      Hours = 7;   Date = (Get-Date -Year 2012 -Month 2 -Day 19 -Hour 0 -Minute 0 -Second 0) },
    @{ Row = 48; Text = "Showing filled tests list and window";                      Hours = 2.5; Date = (Get-Date -Year 2012 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0) }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Text
    $ws.Cells.Item($r.Row, 3).Value = $r.Hours

    # Reuse the existing date-column style (same as D4:D44) instead of letting
    # auto-date-detection mint a brand new number format.
    $ws.Range("D44").Copy()
    $cellD = $ws.Cells.Item($r.Row, 4)
    $cellD.PasteSpecial(-4122)  # xlPasteFormats
    $cellD.Value = $r.Date
}

# Update the view state to match the target workbook
$window = $ws.Application.ActiveWindow
$window.ScrollRow = 20
$ws.Range("B40").Select()
